$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A12 must hold the text "05.07.2019" as a shared string, using the same
# wrap-text style as the other date cells in column A (style index 6,
# e.g. the style already used by A10). Setting a date-shaped string
# directly makes Excel auto-convert it to a date serial, so first coerce
# the cell to Text format (reusing the existing Text style from B2,
# rather than minting a new numFmt/style) before typing the value, then
# restore the correct wrap-text date-column style by copying it from A10.
$ws.Range("B2").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "05.07.2019"
$ws.Range("A10").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B12").Value = "Added a component to correctly display cards' manacosts"
$ws.Range("C12").Value = 2

$ws.Range("B12").Select() | Out-Null
